# Scheduled runner refresh: update cached market-board price/profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the per-job Sheets.
# Values below were re-pulled from the pricing source; a few rows also
# gain/lose their HQ-profit (N) cell depending on whether an HQ price now
# exists for that leve.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 978.6
$ws.Range("I19").Value = 966
$ws.Range("K19").Value = 966
$ws.Range("M19").Value = -791
$ws.Range("H42").Value = 481.6
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H62").Value = 7417.2
$ws.Range("I62").Value = 8521.5
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 8521.5
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -7897.5
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 7417.2
$ws.Range("I65").Value = 8521.5
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 42607.5
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -39487.5
$ws.Range("N65").Value = -21240
$ws.Range("H80").Value = 1940.3846
$ws.Range("I80").Value = 2985.5715
$ws.Range("J80").Value = 721
$ws.Range("K80").Value = 8956.7145
$ws.Range("L80").Value = 2163
$ws.Range("M80").Value = -7958.7145
$ws.Range("N80").Value = -4159
$ws.Range("H83").Value = 1940.3846
$ws.Range("I83").Value = 2985.5715
$ws.Range("J83").Value = 721
$ws.Range("K83").Value = 26870.1435
$ws.Range("L83").Value = 6489
$ws.Range("M83").Value = -21878.1435
$ws.Range("N83").Value = -16473
$ws.Range("H88").Value = 1748.75
$ws.Range("I88").Value = 1497.5
$ws.Range("J88").Value = 2000
$ws.Range("K88").Value = 1497.5
$ws.Range("L88").Value = 2000
$ws.Range("M88").Value = -1091.5
$ws.Range("N88").Value = -2812
$ws.Range("H91").Value = 1748.75
$ws.Range("I91").Value = 1497.5
$ws.Range("J91").Value = 2000
$ws.Range("K91").Value = 1497.5
$ws.Range("L91").Value = 2000
$ws.Range("M91").Value = -93.5
$ws.Range("N91").Value = -4808
$ws.Range("H92").Value = 718.9524
$ws.Range("I92").Value = 783.4375
$ws.Range("K92").Value = 783.4375
$ws.Range("M92").Value = 464.5625
$ws.Range("H116").Value = 4497.5
$ws.Range("I116").Value = 4497
$ws.Range("K116").Value = 4497
$ws.Range("M116").Value = -1055
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2481.6667
$ws.Range("I2").Value = 2481.6667
$ws.Range("K2").Value = 2481.6667
$ws.Range("M2").Value = -2368.6667
$ws.Range("H32").Value = 6315.087
$ws.Range("I32").Value = 6315.087
$ws.Range("K32").Value = 6315.087
$ws.Range("M32").Value = -6028.087
$ws.Range("H61").Value = 1789.2222
$ws.Range("I61").Value = 1789.2222
$ws.Range("K61").Value = 1789.2222
$ws.Range("M61").Value = -1577.2222
$ws.Range("H110").Value = 2673.75
$ws.Range("I110").Value = 2870
$ws.Range("J110").Value = 1300
$ws.Range("K110").Value = 2870
$ws.Range("L110").Value = 1300
$ws.Range("M110").Value = -825
$ws.Range("N110").Value = -5390
$ws.Range("H116").Value = 2481.6667
$ws.Range("I116").Value = 2481.6667
$ws.Range("K116").Value = 2481.6667
$ws.Range("M116").Value = -187.6667000000002
$ws.Range("H132").Value = 5361.5
$ws.Range("I132").Value = 5966.6665
$ws.Range("J132").Value = 4998.4
$ws.Range("K132").Value = 17899.9995
$ws.Range("L132").Value = 14995.2
$ws.Range("M132").Value = -15369.9995
$ws.Range("N132").Value = -20055.2
$ws.Range("H136").Value = 1789.2222
$ws.Range("I136").Value = 1789.2222
$ws.Range("K136").Value = 5367.6666
$ws.Range("M136").Value = -2817.6666
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2481.6667
$ws.Range("I3").Value = 2481.6667
$ws.Range("K3").Value = 2481.6667
$ws.Range("M3").Value = -2367.6667
$ws.Range("H86").Value = 7223.143
$ws.Range("I86").Value = 4333
$ws.Range("J86").Value = 8379.200000000001
$ws.Range("K86").Value = 4333
$ws.Range("L86").Value = 8379.200000000001
$ws.Range("M86").Value = -3210
$ws.Range("N86").Value = -10625.2
$ws.Range("H89").Value = 7223.143
$ws.Range("I89").Value = 4333
$ws.Range("J89").Value = 8379.200000000001
$ws.Range("K89").Value = 21665
$ws.Range("L89").Value = 41896
$ws.Range("M89").Value = -16049
$ws.Range("N89").Value = -53128
$ws.Range("H99").Value = 882.1667
$ws.Range("I99").Value = 882.1667
$ws.Range("K99").Value = 882.1667
$ws.Range("M99").Value = 615.8333
$ws.Range("H107").Value = 1924.5
$ws.Range("I107").Value = 1924.5
$ws.Range("K107").Value = 1924.5
$ws.Range("M107").Value = -4.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 998.5
$ws.Range("I16").Value = 998
$ws.Range("K16").Value = 998
$ws.Range("M16").Value = -711
$ws.Range("H31").Value = 3202.2
$ws.Range("I31").Value = 2937.3333
$ws.Range("K31").Value = 2937.3333
$ws.Range("M31").Value = -2642.3333
$ws.Range("H34").Value = 3202.2
$ws.Range("I34").Value = 2937.3333
$ws.Range("K34").Value = 2937.3333
$ws.Range("M34").Value = -2735.3333
$ws.Range("H113").Value = 998.5
$ws.Range("I113").Value = 998
$ws.Range("K113").Value = 998
$ws.Range("M113").Value = 1172
$ws.Range("H132").Value = 2254.1875
$ws.Range("I132").Value = 1759.6666
$ws.Range("K132").Value = 5278.9998
$ws.Range("M132").Value = -2748.9998
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 3250
$ws.Range("J25").Value = 5000
$ws.Range("L25").Value = 15000
$ws.Range("N25").Value = -15338
$ws.Range("H30").Value = 3250
$ws.Range("J30").Value = 5000
$ws.Range("L30").Value = 15000
$ws.Range("N30").Value = -15204
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1855.6666
$ws.Range("I80").Value = 1282.6666
$ws.Range("K80").Value = 1282.6666
$ws.Range("M80").Value = -284.6666
$ws.Range("H83").Value = 1855.6666
$ws.Range("I83").Value = 1282.6666
$ws.Range("K83").Value = 6413.333000000001
$ws.Range("M83").Value = -1421.333000000001
$ws.Range("H97").Value = 3414.8333
$ws.Range("I97").Value = 3414.8333
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 3414.8333
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -2918.8333
$ws.Range("N97").ClearContents()
$ws.Range("H107").Value = 3597.75
$ws.Range("I107").Value = 3597.75
$ws.Range("K107").Value = 3597.75
$ws.Range("M107").Value = -1677.75
$ws.Range("H122").Value = 4499.25
$ws.Range("I122").Value = 4499.25
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13497.75
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11047.75
$ws.Range("N122").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 47500
$ws.Range("J88").Value = 47500
$ws.Range("L88").Value = 47500
$ws.Range("N88").Value = -48356
$ws.Range("H91").Value = 47500
$ws.Range("J91").Value = 47500
$ws.Range("L91").Value = 47500
$ws.Range("N91").Value = -50464
$ws.Range("H93").Value = 3780
$ws.Range("J93").Value = 1000
$ws.Range("L93").Value = 1000
$ws.Range("N93").Value = -3496
$ws.Range("H132").Value = 4617.6816
$ws.Range("I132").Value = 3601.2727
$ws.Range("J132").Value = 5634.091
$ws.Range("K132").Value = 10803.8181
$ws.Range("L132").Value = 16902.273
$ws.Range("M132").Value = -8273.8181
$ws.Range("N132").Value = -21962.273
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1931.5454
$ws.Range("I81").Value = 2241.5
$ws.Range("J81").Value = 1559.6
$ws.Range("K81").Value = 4483
$ws.Range("L81").Value = 3119.2
$ws.Range("M81").Value = -3422
$ws.Range("N81").Value = -5241.2
$ws.Range("H84").Value = 1931.5454
$ws.Range("I84").Value = 2241.5
$ws.Range("J84").Value = 1559.6
$ws.Range("K84").Value = 22415
$ws.Range("L84").Value = 15596
$ws.Range("M84").Value = -17111
$ws.Range("N84").Value = -26204
$ws.Range("H132").Value = 2571
$ws.Range("I132").Value = 1637.3846
$ws.Range("K132").Value = 4912.1538
$ws.Range("M132").Value = -2382.1538
